$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" column header in H1 -- copy formatting (bold, border, centered)
# from the neighboring header cell G1, then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0
$ws.Range("H1").Value = "Save"

# New value for the "Save" column in the data row.
$ws.Range("H2").Value = 1
